# Event.xlsx - "출석이벤트 5→7일 변경 및 보상 추가"
# (Attendance event changed from 5 to 7 days; rewards added)

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Event_Main")
$wsList = $wb.Worksheets.Item("Event_List")

# --- Update existing reward rows on Event_List ---
# Row 4 (132003 / day 3) now grants item 201007 x1
$wsList.Range("D4").Value = 201007
$wsList.Range("E4").Value = 1

# Row 5 (132004 / day 4) now grants item 100001 x3000
$wsList.Range("D5").Value = 100001
$wsList.Range("E5").Value = 3000

# Row 6 (132005 / day 5) now grants item 100002 x200
$wsList.Range("D6").Value = 100002
$wsList.Range("E6").Value = 200

# --- Add new day 6 / day 7 reward rows ---
# Row 7: 132006 / day 6 -> item 331002 x1
$wsList.Range("A7").Value = 132006
$wsList.Range("B7").Value = 32001
$wsList.Range("C7").Value = 6
$wsList.Range("D7").Value = 331002
$wsList.Range("E7").Value = 1
$wsList.Range("F7").Value = 45573
$wsList.Range("G7").Value = 72686
$wsList.Range("H7").Value = $false

# Row 8: 132007 / day 7 -> item 322002 x1
$wsList.Range("A8").Value = 132007
$wsList.Range("B8").Value = 32001
$wsList.Range("C8").Value = 7
$wsList.Range("D8").Value = 322002
$wsList.Range("E8").Value = 1
$wsList.Range("F8").Value = 45573
$wsList.Range("G8").Value = 72686
$wsList.Range("H8").Value = $false

# Copy the row-6 cell formatting down onto the two new rows so the
# new cells (number formats / alignment) match the rest of the table.
$wsList.Range("A6:H6").Copy()
[void]$wsList.Range("A7:H7").PasteSpecial(-4122)
$wsList.Range("A6:H6").Copy()
[void]$wsList.Range("A8:H8").PasteSpecial(-4122)

# --- Selection / active sheet bookkeeping ---
[void]$wsMain.Range("D6").Select()
[void]$wsList.Range("F11").Select()
[void]$wsList.Activate()
